$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dimension implicitly by writing to the full new range (A1:C181).
# Row 1 (header: ano / mes / valor) is unchanged.

$ws.Cells.Item(2, 1).Value = 2010
$ws.Cells.Item(2, 2).Value = "Abril"
$ws.Cells.Item(2, 3).Value = 1.75657
$ws.Cells.Item(3, 1).Value = 2010
$ws.Cells.Item(3, 2).Value = "Agosto"
$ws.Cells.Item(3, 3).Value = 1.759563636363636
$ws.Cells.Item(4, 1).Value = 2010
$ws.Cells.Item(4, 2).Value = "Dezembro"
$ws.Cells.Item(4, 3).Value = 1.693413043478261
$ws.Cells.Item(5, 1).Value = 2010
$ws.Cells.Item(5, 2).Value = "Fevereiro"
$ws.Cells.Item(5, 3).Value = 1.841633333333333
$ws.Cells.Item(6, 1).Value = 2010
$ws.Cells.Item(6, 2).Value = "Janeiro"
$ws.Cells.Item(6, 3).Value = 1.77982
$ws.Cells.Item(7, 1).Value = 2010
$ws.Cells.Item(7, 2).Value = "Julho"
$ws.Cells.Item(7, 3).Value = 1.769636363636364
$ws.Cells.Item(8, 1).Value = 2010
$ws.Cells.Item(8, 2).Value = "Junho"
$ws.Cells.Item(8, 3).Value = 1.806528571428571
$ws.Cells.Item(9, 1).Value = 2010
$ws.Cells.Item(9, 2).Value = "Maio"
$ws.Cells.Item(9, 3).Value = 1.813190476190476
$ws.Cells.Item(10, 1).Value = 2010
$ws.Cells.Item(10, 2).Value = "Março"
$ws.Cells.Item(10, 3).Value = 1.785843478260869
$ws.Cells.Item(11, 1).Value = 2010
$ws.Cells.Item(11, 2).Value = "Novembro"
$ws.Cells.Item(11, 3).Value = 1.71333
$ws.Cells.Item(12, 1).Value = 2010
$ws.Cells.Item(12, 2).Value = "Outubro"
$ws.Cells.Item(12, 3).Value = 1.6835
$ws.Cells.Item(13, 1).Value = 2010
$ws.Cells.Item(13, 2).Value = "Setembro"
$ws.Cells.Item(13, 3).Value = 1.718709523809524
$ws.Cells.Item(14, 1).Value = 2011
$ws.Cells.Item(14, 2).Value = "Abril"
$ws.Cells.Item(14, 3).Value = 1.586447368421052
$ws.Cells.Item(15, 1).Value = 2011
$ws.Cells.Item(15, 2).Value = "Agosto"
$ws.Cells.Item(15, 3).Value = 1.597008695652174
$ws.Cells.Item(16, 1).Value = 2011
$ws.Cells.Item(16, 2).Value = "Dezembro"
$ws.Cells.Item(16, 3).Value = 1.836886363636363
$ws.Cells.Item(17, 1).Value = 2011
$ws.Cells.Item(17, 2).Value = "Fevereiro"
$ws.Cells.Item(17, 3).Value = 1.66799
$ws.Cells.Item(18, 1).Value = 2011
$ws.Cells.Item(18, 2).Value = "Janeiro"
$ws.Cells.Item(18, 3).Value = 1.674914285714286
$ws.Cells.Item(19, 1).Value = 2011
$ws.Cells.Item(19, 2).Value = "Julho"
$ws.Cells.Item(19, 3).Value = 1.563938095238095
$ws.Cells.Item(20, 1).Value = 2011
$ws.Cells.Item(20, 2).Value = "Junho"
$ws.Cells.Item(20, 3).Value = 1.587042857142857
$ws.Cells.Item(21, 1).Value = 2011
$ws.Cells.Item(21, 2).Value = "Maio"
$ws.Cells.Item(21, 3).Value = 1.613490909090909
$ws.Cells.Item(22, 1).Value = 2011
$ws.Cells.Item(22, 2).Value = "Março"
$ws.Cells.Item(22, 3).Value = 1.6591
$ws.Cells.Item(23, 1).Value = 2011
$ws.Cells.Item(23, 2).Value = "Novembro"
$ws.Cells.Item(23, 3).Value = 1.79049
$ws.Cells.Item(24, 1).Value = 2011
$ws.Cells.Item(24, 2).Value = "Outubro"
$ws.Cells.Item(24, 3).Value = 1.77257
$ws.Cells.Item(25, 1).Value = 2011
$ws.Cells.Item(25, 2).Value = "Setembro"
$ws.Cells.Item(25, 3).Value = 1.749776190476191
$ws.Cells.Item(26, 1).Value = 2012
$ws.Cells.Item(26, 2).Value = "Abril"
$ws.Cells.Item(26, 3).Value = 1.854835
$ws.Cells.Item(27, 1).Value = 2012
$ws.Cells.Item(27, 2).Value = "Agosto"
$ws.Cells.Item(27, 3).Value = 2.02944347826087
$ws.Cells.Item(28, 1).Value = 2012
$ws.Cells.Item(28, 2).Value = "Dezembro"
$ws.Cells.Item(28, 3).Value = 2.077835
$ws.Cells.Item(29, 1).Value = 2012
$ws.Cells.Item(29, 2).Value = "Fevereiro"
$ws.Cells.Item(29, 3).Value = 1.718394736842106
$ws.Cells.Item(30, 1).Value = 2012
$ws.Cells.Item(30, 2).Value = "Janeiro"
$ws.Cells.Item(30, 3).Value = 1.789681818181818
$ws.Cells.Item(31, 1).Value = 2012
$ws.Cells.Item(31, 2).Value = "Julho"
$ws.Cells.Item(31, 3).Value = 2.028736363636364
$ws.Cells.Item(32, 1).Value = 2012
$ws.Cells.Item(32, 2).Value = "Junho"
$ws.Cells.Item(32, 3).Value = 2.049195
$ws.Cells.Item(33, 1).Value = 2012
$ws.Cells.Item(33, 2).Value = "Maio"
$ws.Cells.Item(33, 3).Value = 1.985990909090909
$ws.Cells.Item(34, 1).Value = 2012
$ws.Cells.Item(34, 2).Value = "Março"
$ws.Cells.Item(34, 3).Value = 1.795309090909091
$ws.Cells.Item(35, 1).Value = 2012
$ws.Cells.Item(35, 2).Value = "Novembro"
$ws.Cells.Item(35, 3).Value = 2.06775
$ws.Cells.Item(36, 1).Value = 2012
$ws.Cells.Item(36, 2).Value = "Outubro"
$ws.Cells.Item(36, 3).Value = 2.029845454545455
$ws.Cells.Item(37, 1).Value = 2012
$ws.Cells.Item(37, 2).Value = "Setembro"
$ws.Cells.Item(37, 3).Value = 2.028078947368421
$ws.Cells.Item(38, 1).Value = 2013
$ws.Cells.Item(38, 2).Value = "Abril"
$ws.Cells.Item(38, 3).Value = 2.002213636363636
$ws.Cells.Item(39, 1).Value = 2013
$ws.Cells.Item(39, 2).Value = "Agosto"
$ws.Cells.Item(39, 3).Value = 2.342190909090909
$ws.Cells.Item(40, 1).Value = 2013
$ws.Cells.Item(40, 2).Value = "Dezembro"
$ws.Cells.Item(40, 3).Value = 2.345485714285714
$ws.Cells.Item(41, 1).Value = 2013
$ws.Cells.Item(41, 2).Value = "Fevereiro"
$ws.Cells.Item(41, 3).Value = 1.97325
$ws.Cells.Item(42, 1).Value = 2013
$ws.Cells.Item(42, 2).Value = "Janeiro"
$ws.Cells.Item(42, 3).Value = 2.031077272727273
$ws.Cells.Item(43, 1).Value = 2013
$ws.Cells.Item(43, 2).Value = "Julho"
$ws.Cells.Item(43, 3).Value = 2.252169565217391
$ws.Cells.Item(44, 1).Value = 2013
$ws.Cells.Item(44, 2).Value = "Junho"
$ws.Cells.Item(44, 3).Value = 2.172955
$ws.Cells.Item(45, 1).Value = 2013
$ws.Cells.Item(45, 2).Value = "Maio"
$ws.Cells.Item(45, 3).Value = 2.034842857142857
$ws.Cells.Item(46, 1).Value = 2013
$ws.Cells.Item(46, 2).Value = "Março"
$ws.Cells.Item(46, 3).Value = 1.98284
$ws.Cells.Item(47, 1).Value = 2013
$ws.Cells.Item(47, 2).Value = "Novembro"
$ws.Cells.Item(47, 3).Value = 2.29535
$ws.Cells.Item(48, 1).Value = 2013
$ws.Cells.Item(48, 2).Value = "Outubro"
$ws.Cells.Item(48, 3).Value = 2.188647826086957
$ws.Cells.Item(49, 1).Value = 2013
$ws.Cells.Item(49, 2).Value = "Setembro"
$ws.Cells.Item(49, 3).Value = 2.270509523809524
$ws.Cells.Item(50, 1).Value = 2014
$ws.Cells.Item(50, 2).Value = "Abril"
$ws.Cells.Item(50, 3).Value = 2.23277
$ws.Cells.Item(51, 1).Value = 2014
$ws.Cells.Item(51, 2).Value = "Agosto"
$ws.Cells.Item(51, 3).Value = 2.268028571428571
$ws.Cells.Item(52, 1).Value = 2014
$ws.Cells.Item(52, 2).Value = "Dezembro"
$ws.Cells.Item(52, 3).Value = 2.639363636363636
$ws.Cells.Item(53, 1).Value = 2014
$ws.Cells.Item(53, 2).Value = "Fevereiro"
$ws.Cells.Item(53, 3).Value = 2.38368
$ws.Cells.Item(54, 1).Value = 2014
$ws.Cells.Item(54, 2).Value = "Janeiro"
$ws.Cells.Item(54, 3).Value = 2.382209090909091
$ws.Cells.Item(55, 1).Value = 2014
$ws.Cells.Item(55, 2).Value = "Julho"
$ws.Cells.Item(55, 3).Value = 2.224647826086957
$ws.Cells.Item(56, 1).Value = 2014
$ws.Cells.Item(56, 2).Value = "Junho"
$ws.Cells.Item(56, 3).Value = 2.23547
$ws.Cells.Item(57, 1).Value = 2014
$ws.Cells.Item(57, 2).Value = "Maio"
$ws.Cells.Item(57, 3).Value = 2.220880952380953
$ws.Cells.Item(58, 1).Value = 2014
$ws.Cells.Item(58, 2).Value = "Março"
$ws.Cells.Item(58, 3).Value = 2.326089473684211
$ws.Cells.Item(59, 1).Value = 2014
$ws.Cells.Item(59, 2).Value = "Novembro"
$ws.Cells.Item(59, 3).Value = 2.548365
$ws.Cells.Item(60, 1).Value = 2014
$ws.Cells.Item(60, 2).Value = "Outubro"
$ws.Cells.Item(60, 3).Value = 2.448260869565217
$ws.Cells.Item(61, 1).Value = 2014
$ws.Cells.Item(61, 2).Value = "Setembro"
$ws.Cells.Item(61, 3).Value = 2.332868181818182
$ws.Cells.Item(62, 1).Value = 2015
$ws.Cells.Item(62, 2).Value = "Abril"
$ws.Cells.Item(62, 3).Value = 3.04322
$ws.Cells.Item(63, 1).Value = 2015
$ws.Cells.Item(63, 2).Value = "Agosto"
$ws.Cells.Item(63, 3).Value = 3.514304761904762
$ws.Cells.Item(64, 1).Value = 2015
$ws.Cells.Item(64, 2).Value = "Dezembro"
$ws.Cells.Item(64, 3).Value = 3.871136363636363
$ws.Cells.Item(65, 1).Value = 2015
$ws.Cells.Item(65, 2).Value = "Fevereiro"
$ws.Cells.Item(65, 3).Value = 2.81645
$ws.Cells.Item(66, 1).Value = 2015
$ws.Cells.Item(66, 2).Value = "Janeiro"
$ws.Cells.Item(66, 3).Value = 2.634228571428571
$ws.Cells.Item(67, 1).Value = 2015
$ws.Cells.Item(67, 2).Value = "Julho"
$ws.Cells.Item(67, 3).Value = 3.223143478260869
$ws.Cells.Item(68, 1).Value = 2015
$ws.Cells.Item(68, 2).Value = "Junho"
$ws.Cells.Item(68, 3).Value = 3.111738095238096
$ws.Cells.Item(69, 1).Value = 2015
$ws.Cells.Item(69, 2).Value = "Maio"
$ws.Cells.Item(69, 3).Value = 3.061715
$ws.Cells.Item(70, 1).Value = 2015
$ws.Cells.Item(70, 2).Value = "Março"
$ws.Cells.Item(70, 3).Value = 3.139477272727273
$ws.Cells.Item(71, 1).Value = 2015
$ws.Cells.Item(71, 2).Value = "Novembro"
$ws.Cells.Item(71, 3).Value = 3.77646
$ws.Cells.Item(72, 1).Value = 2015
$ws.Cells.Item(72, 2).Value = "Outubro"
$ws.Cells.Item(72, 3).Value = 3.880138095238095
$ws.Cells.Item(73, 1).Value = 2015
$ws.Cells.Item(73, 2).Value = "Setembro"
$ws.Cells.Item(73, 3).Value = 3.906457142857143
$ws.Cells.Item(74, 1).Value = 2016
$ws.Cells.Item(74, 2).Value = "Abril"
$ws.Cells.Item(74, 3).Value = 3.565845
$ws.Cells.Item(75, 1).Value = 2016
$ws.Cells.Item(75, 2).Value = "Agosto"
$ws.Cells.Item(75, 3).Value = 3.209660869565217
$ws.Cells.Item(76, 1).Value = 2016
$ws.Cells.Item(76, 2).Value = "Dezembro"
$ws.Cells.Item(76, 3).Value = 3.352268181818182
$ws.Cells.Item(77, 1).Value = 2016
$ws.Cells.Item(77, 2).Value = "Fevereiro"
$ws.Cells.Item(77, 3).Value = 3.973742105263157
$ws.Cells.Item(78, 1).Value = 2016
$ws.Cells.Item(78, 2).Value = "Janeiro"
$ws.Cells.Item(78, 3).Value = 4.05235
$ws.Cells.Item(79, 1).Value = 2016
$ws.Cells.Item(79, 2).Value = "Julho"
$ws.Cells.Item(79, 3).Value = 3.275566666666667
$ws.Cells.Item(80, 1).Value = 2016
$ws.Cells.Item(80, 2).Value = "Junho"
$ws.Cells.Item(80, 3).Value = 3.424477272727273
$ws.Cells.Item(81, 1).Value = 2016
$ws.Cells.Item(81, 2).Value = "Maio"
$ws.Cells.Item(81, 3).Value = 3.539290476190476
$ws.Cells.Item(82, 1).Value = 2016
$ws.Cells.Item(82, 2).Value = "Março"
$ws.Cells.Item(82, 3).Value = 3.703918181818182
$ws.Cells.Item(83, 1).Value = 2016
$ws.Cells.Item(83, 2).Value = "Novembro"
$ws.Cells.Item(83, 3).Value = 3.34203
$ws.Cells.Item(84, 1).Value = 2016
$ws.Cells.Item(84, 2).Value = "Outubro"
$ws.Cells.Item(84, 3).Value = 3.185845
$ws.Cells.Item(85, 1).Value = 2016
$ws.Cells.Item(85, 2).Value = "Setembro"
$ws.Cells.Item(85, 3).Value = 3.256371428571428
$ws.Cells.Item(86, 1).Value = 2017
$ws.Cells.Item(86, 2).Value = "Abril"
$ws.Cells.Item(86, 3).Value = 3.136172222222222
$ws.Cells.Item(87, 1).Value = 2017
$ws.Cells.Item(87, 2).Value = "Agosto"
$ws.Cells.Item(87, 3).Value = 3.150917391304348
$ws.Cells.Item(88, 1).Value = 2017
$ws.Cells.Item(88, 2).Value = "Dezembro"
$ws.Cells.Item(88, 3).Value = 3.291915
$ws.Cells.Item(89, 1).Value = 2017
$ws.Cells.Item(89, 2).Value = "Fevereiro"
$ws.Cells.Item(89, 3).Value = 3.104194444444444
$ws.Cells.Item(90, 1).Value = 2017
$ws.Cells.Item(90, 2).Value = "Janeiro"
$ws.Cells.Item(90, 3).Value = 3.196609090909091
$ws.Cells.Item(91, 1).Value = 2017
$ws.Cells.Item(91, 2).Value = "Julho"
$ws.Cells.Item(91, 3).Value = 3.206138095238095
$ws.Cells.Item(92, 1).Value = 2017
$ws.Cells.Item(92, 2).Value = "Junho"
$ws.Cells.Item(92, 3).Value = 3.295366666666666
$ws.Cells.Item(93, 1).Value = 2017
$ws.Cells.Item(93, 2).Value = "Maio"
$ws.Cells.Item(93, 3).Value = 3.209509090909091
$ws.Cells.Item(94, 1).Value = 2017
$ws.Cells.Item(94, 2).Value = "Março"
$ws.Cells.Item(94, 3).Value = 3.127930434782609
$ws.Cells.Item(95, 1).Value = 2017
$ws.Cells.Item(95, 2).Value = "Novembro"
$ws.Cells.Item(95, 3).Value = 3.25938
$ws.Cells.Item(96, 1).Value = 2017
$ws.Cells.Item(96, 2).Value = "Outubro"
$ws.Cells.Item(96, 3).Value = 3.191228571428571
$ws.Cells.Item(97, 1).Value = 2017
$ws.Cells.Item(97, 2).Value = "Setembro"
$ws.Cells.Item(97, 3).Value = 3.13479
$ws.Cells.Item(98, 1).Value = 2018
$ws.Cells.Item(98, 2).Value = "Abril"
$ws.Cells.Item(98, 3).Value = 3.407495238095238
$ws.Cells.Item(99, 1).Value = 2018
$ws.Cells.Item(99, 2).Value = "Agosto"
$ws.Cells.Item(99, 3).Value = 3.92975652173913
$ws.Cells.Item(100, 1).Value = 2018
$ws.Cells.Item(100, 2).Value = "Dezembro"
$ws.Cells.Item(100, 3).Value = 3.885055
$ws.Cells.Item(101, 1).Value = 2018
$ws.Cells.Item(101, 2).Value = "Fevereiro"
$ws.Cells.Item(101, 3).Value = 3.2415
$ws.Cells.Item(102, 1).Value = 2018
$ws.Cells.Item(102, 2).Value = "Janeiro"
$ws.Cells.Item(102, 3).Value = 3.210609090909091
$ws.Cells.Item(103, 1).Value = 2018
$ws.Cells.Item(103, 2).Value = "Julho"
$ws.Cells.Item(103, 3).Value = 3.828763636363636
$ws.Cells.Item(104, 1).Value = 2018
$ws.Cells.Item(104, 2).Value = "Junho"
$ws.Cells.Item(104, 3).Value = 3.773171428571428
$ws.Cells.Item(105, 1).Value = 2018
$ws.Cells.Item(105, 2).Value = "Maio"
$ws.Cells.Item(105, 3).Value = 3.636057142857143
$ws.Cells.Item(106, 1).Value = 2018
$ws.Cells.Item(106, 2).Value = "Março"
$ws.Cells.Item(106, 3).Value = 3.279214285714286
$ws.Cells.Item(107, 1).Value = 2018
$ws.Cells.Item(107, 2).Value = "Novembro"
$ws.Cells.Item(107, 3).Value = 3.786665
$ws.Cells.Item(108, 1).Value = 2018
$ws.Cells.Item(108, 2).Value = "Outubro"
$ws.Cells.Item(108, 3).Value = 3.758409090909091
$ws.Cells.Item(109, 1).Value = 2018
$ws.Cells.Item(109, 2).Value = "Setembro"
$ws.Cells.Item(109, 3).Value = 4.116547368421053
$ws.Cells.Item(110, 1).Value = 2019
$ws.Cells.Item(110, 2).Value = "Abril"
$ws.Cells.Item(110, 3).Value = 3.896157142857143
$ws.Cells.Item(111, 1).Value = 2019
$ws.Cells.Item(111, 2).Value = "Agosto"
$ws.Cells.Item(111, 3).Value = 4.019981818181818
$ws.Cells.Item(112, 1).Value = 2019
$ws.Cells.Item(112, 2).Value = "Dezembro"
$ws.Cells.Item(112, 3).Value = 4.109590476190476
$ws.Cells.Item(113, 1).Value = 2019
$ws.Cells.Item(113, 2).Value = "Fevereiro"
$ws.Cells.Item(113, 3).Value = 3.723625
$ws.Cells.Item(114, 1).Value = 2019
$ws.Cells.Item(114, 2).Value = "Janeiro"
$ws.Cells.Item(114, 3).Value = 3.741681818181819
$ws.Cells.Item(115, 1).Value = 2019
$ws.Cells.Item(115, 2).Value = "Julho"
$ws.Cells.Item(115, 3).Value = 3.779339130434783
$ws.Cells.Item(116, 1).Value = 2019
$ws.Cells.Item(116, 2).Value = "Junho"
$ws.Cells.Item(116, 3).Value = 3.858826315789474
$ws.Cells.Item(117, 1).Value = 2019
$ws.Cells.Item(117, 2).Value = "Maio"
$ws.Cells.Item(117, 3).Value = 4.001518181818182
$ws.Cells.Item(118, 1).Value = 2019
$ws.Cells.Item(118, 2).Value = "Março"
$ws.Cells.Item(118, 3).Value = 3.846484210526316
$ws.Cells.Item(119, 1).Value = 2019
$ws.Cells.Item(119, 2).Value = "Novembro"
$ws.Cells.Item(119, 3).Value = 4.155345
$ws.Cells.Item(120, 1).Value = 2019
$ws.Cells.Item(120, 2).Value = "Outubro"
$ws.Cells.Item(120, 3).Value = 4.086986956521739
$ws.Cells.Item(121, 1).Value = 2019
$ws.Cells.Item(121, 2).Value = "Setembro"
$ws.Cells.Item(121, 3).Value = 4.1215
$ws.Cells.Item(122, 1).Value = 2020
$ws.Cells.Item(122, 2).Value = "Abril"
$ws.Cells.Item(122, 3).Value = 5.32558
$ws.Cells.Item(123, 1).Value = 2020
$ws.Cells.Item(123, 2).Value = "Agosto"
$ws.Cells.Item(123, 3).Value = 5.461233333333333
$ws.Cells.Item(124, 1).Value = 2020
$ws.Cells.Item(124, 2).Value = "Dezembro"
$ws.Cells.Item(124, 3).Value = 5.145586363636363
$ws.Cells.Item(125, 1).Value = 2020
$ws.Cells.Item(125, 2).Value = "Fevereiro"
$ws.Cells.Item(125, 3).Value = 4.341011111111111
$ws.Cells.Item(126, 1).Value = 2020
$ws.Cells.Item(126, 2).Value = "Janeiro"
$ws.Cells.Item(126, 3).Value = 4.149463636363636
$ws.Cells.Item(127, 1).Value = 2020
$ws.Cells.Item(127, 2).Value = "Julho"
$ws.Cells.Item(127, 3).Value = 5.280191304347826
$ws.Cells.Item(128, 1).Value = 2020
$ws.Cells.Item(128, 2).Value = "Junho"
$ws.Cells.Item(128, 3).Value = 5.1966
$ws.Cells.Item(129, 1).Value = 2020
$ws.Cells.Item(129, 2).Value = "Maio"
$ws.Cells.Item(129, 3).Value = 5.643445
$ws.Cells.Item(130, 1).Value = 2020
$ws.Cells.Item(130, 2).Value = "Março"
$ws.Cells.Item(130, 3).Value = 4.883854545454546
$ws.Cells.Item(131, 1).Value = 2020
$ws.Cells.Item(131, 2).Value = "Novembro"
$ws.Cells.Item(131, 3).Value = 5.417835
$ws.Cells.Item(132, 1).Value = 2020
$ws.Cells.Item(132, 2).Value = "Outubro"
$ws.Cells.Item(132, 3).Value = 5.625790476190476
$ws.Cells.Item(133, 1).Value = 2020
$ws.Cells.Item(133, 2).Value = "Setembro"
$ws.Cells.Item(133, 3).Value = 5.399485714285714
$ws.Cells.Item(134, 1).Value = 2021
$ws.Cells.Item(134, 2).Value = "Abril"
$ws.Cells.Item(134, 3).Value = 5.562135
$ws.Cells.Item(135, 1).Value = 2021
$ws.Cells.Item(135, 2).Value = "Agosto"
$ws.Cells.Item(135, 3).Value = 5.251718181818182
$ws.Cells.Item(136, 1).Value = 2021
$ws.Cells.Item(136, 2).Value = "Dezembro"
$ws.Cells.Item(136, 3).Value = 5.651391304347826
$ws.Cells.Item(137, 1).Value = 2021
$ws.Cells.Item(137, 2).Value = "Fevereiro"
$ws.Cells.Item(137, 3).Value = 5.416494444444444
$ws.Cells.Item(138, 1).Value = 2021
$ws.Cells.Item(138, 2).Value = "Janeiro"
$ws.Cells.Item(138, 3).Value = 5.356244999999999
$ws.Cells.Item(139, 1).Value = 2021
$ws.Cells.Item(139, 2).Value = "Julho"
$ws.Cells.Item(139, 3).Value = 5.156704545454546
$ws.Cells.Item(140, 1).Value = 2021
$ws.Cells.Item(140, 2).Value = "Junho"
$ws.Cells.Item(140, 3).Value = 5.031904761904762
$ws.Cells.Item(141, 1).Value = 2021
$ws.Cells.Item(141, 2).Value = "Maio"
$ws.Cells.Item(141, 3).Value = 5.291057142857143
$ws.Cells.Item(142, 1).Value = 2021
$ws.Cells.Item(142, 2).Value = "Março"
$ws.Cells.Item(142, 3).Value = 5.646147826086957
$ws.Cells.Item(143, 1).Value = 2021
$ws.Cells.Item(143, 2).Value = "Novembro"
$ws.Cells.Item(143, 3).Value = 5.556859999999999
$ws.Cells.Item(144, 1).Value = 2021
$ws.Cells.Item(144, 2).Value = "Outubro"
$ws.Cells.Item(144, 3).Value = 5.53998
$ws.Cells.Item(145, 1).Value = 2021
$ws.Cells.Item(145, 2).Value = "Setembro"
$ws.Cells.Item(145, 3).Value = 5.279690476190477
$ws.Cells.Item(146, 1).Value = 2022
$ws.Cells.Item(146, 2).Value = "Abril"
$ws.Cells.Item(146, 3).Value = 4.758015789473684
$ws.Cells.Item(147, 1).Value = 2022
$ws.Cells.Item(147, 2).Value = "Agosto"
$ws.Cells.Item(147, 3).Value = 5.143286956521739
$ws.Cells.Item(148, 1).Value = 2022
$ws.Cells.Item(148, 2).Value = "Dezembro"
$ws.Cells.Item(148, 3).Value = 5.242431818181818
$ws.Cells.Item(149, 1).Value = 2022
$ws.Cells.Item(149, 2).Value = "Fevereiro"
$ws.Cells.Item(149, 3).Value = 5.196578947368421
$ws.Cells.Item(150, 1).Value = 2022
$ws.Cells.Item(150, 2).Value = "Janeiro"
$ws.Cells.Item(150, 3).Value = 5.534104761904762
$ws.Cells.Item(151, 1).Value = 2022
$ws.Cells.Item(151, 2).Value = "Julho"
$ws.Cells.Item(151, 3).Value = 5.368071428571429
$ws.Cells.Item(152, 1).Value = 2022
$ws.Cells.Item(152, 2).Value = "Junho"
$ws.Cells.Item(152, 3).Value = 5.049209523809524
$ws.Cells.Item(153, 1).Value = 2022
$ws.Cells.Item(153, 2).Value = "Maio"
$ws.Cells.Item(153, 3).Value = 4.95505
$ws.Cells.Item(154, 1).Value = 2022
$ws.Cells.Item(154, 2).Value = "Março"
$ws.Cells.Item(154, 3).Value = 4.968381818181818
$ws.Cells.Item(155, 1).Value = 2022
$ws.Cells.Item(155, 2).Value = "Novembro"
$ws.Cells.Item(155, 3).Value = 5.274649999999999
$ws.Cells.Item(156, 1).Value = 2022
$ws.Cells.Item(156, 2).Value = "Outubro"
$ws.Cells.Item(156, 3).Value = 5.2503
$ws.Cells.Item(157, 1).Value = 2022
$ws.Cells.Item(157, 2).Value = "Setembro"
$ws.Cells.Item(157, 3).Value = 5.236957142857143
$ws.Cells.Item(158, 1).Value = 2023
$ws.Cells.Item(158, 2).Value = "Abril"
$ws.Cells.Item(158, 3).Value = 5.019733333333333
$ws.Cells.Item(159, 1).Value = 2023
$ws.Cells.Item(159, 2).Value = "Agosto"
$ws.Cells.Item(159, 3).Value = 4.903543478260869
$ws.Cells.Item(160, 1).Value = 2023
$ws.Cells.Item(160, 2).Value = "Dezembro"
$ws.Cells.Item(160, 3).Value = 4.897245
$ws.Cells.Item(161, 1).Value = 2023
$ws.Cells.Item(161, 2).Value = "Fevereiro"
$ws.Cells.Item(161, 3).Value = 5.171688888888889
$ws.Cells.Item(162, 1).Value = 2023
$ws.Cells.Item(162, 2).Value = "Janeiro"
$ws.Cells.Item(162, 3).Value = 5.200681818181819
$ws.Cells.Item(163, 1).Value = 2023
$ws.Cells.Item(163, 2).Value = "Julho"
$ws.Cells.Item(163, 3).Value = 4.800833333333333
$ws.Cells.Item(164, 1).Value = 2023
$ws.Cells.Item(164, 2).Value = "Junho"
$ws.Cells.Item(164, 3).Value = 4.851566666666667
$ws.Cells.Item(165, 1).Value = 2023
$ws.Cells.Item(165, 2).Value = "Maio"
$ws.Cells.Item(165, 3).Value = 4.982836363636364
$ws.Cells.Item(166, 1).Value = 2023
$ws.Cells.Item(166, 2).Value = "Março"
$ws.Cells.Item(166, 3).Value = 5.211460869565218
$ws.Cells.Item(167, 1).Value = 2023
$ws.Cells.Item(167, 2).Value = "Novembro"
$ws.Cells.Item(167, 3).Value = 4.89834
$ws.Cells.Item(168, 1).Value = 2023
$ws.Cells.Item(168, 2).Value = "Outubro"
$ws.Cells.Item(168, 3).Value = 5.064842857142857
$ws.Cells.Item(169, 1).Value = 2023
$ws.Cells.Item(169, 2).Value = "Setembro"
$ws.Cells.Item(169, 3).Value = 4.93699
$ws.Cells.Item(170, 1).Value = 2024
$ws.Cells.Item(170, 2).Value = "Abril"
$ws.Cells.Item(170, 3).Value = 5.129095454545455
$ws.Cells.Item(171, 1).Value = 2024
$ws.Cells.Item(171, 2).Value = "Agosto"
$ws.Cells.Item(171, 3).Value = 5.552613636363636
$ws.Cells.Item(172, 1).Value = 2024
$ws.Cells.Item(172, 2).Value = "Dezembro"
$ws.Cells.Item(172, 3).Value = 6.097028571428571
$ws.Cells.Item(173, 1).Value = 2024
$ws.Cells.Item(173, 2).Value = "Fevereiro"
$ws.Cells.Item(173, 3).Value = 4.964389473684211
$ws.Cells.Item(174, 1).Value = 2024
$ws.Cells.Item(174, 2).Value = "Janeiro"
$ws.Cells.Item(174, 3).Value = 4.914395454545454
$ws.Cells.Item(175, 1).Value = 2024
$ws.Cells.Item(175, 2).Value = "Julho"
$ws.Cells.Item(175, 3).Value = 5.542047826086956
$ws.Cells.Item(176, 1).Value = 2024
$ws.Cells.Item(176, 2).Value = "Junho"
$ws.Cells.Item(176, 3).Value = 5.388975
$ws.Cells.Item(177, 1).Value = 2024
$ws.Cells.Item(177, 2).Value = "Maio"
$ws.Cells.Item(177, 3).Value = 5.133047619047619
$ws.Cells.Item(178, 1).Value = 2024
$ws.Cells.Item(178, 2).Value = "Março"
$ws.Cells.Item(178, 3).Value = 4.980135
$ws.Cells.Item(179, 1).Value = 2024
$ws.Cells.Item(179, 2).Value = "Novembro"
$ws.Cells.Item(179, 3).Value = 5.807057894736842
$ws.Cells.Item(180, 1).Value = 2024
$ws.Cells.Item(180, 2).Value = "Outubro"
$ws.Cells.Item(180, 3).Value = 5.624108695652174
$ws.Cells.Item(181, 1).Value = 2024
$ws.Cells.Item(181, 2).Value = "Setembro"
$ws.Cells.Item(181, 3).Value = 5.541566666666666

Write-Host "done"